# Apply the shortage-report update: refresh the "balance" ratios for a
# couple of items and bump the generated-at timestamp shown at the
# bottom of the sheet (re-uploaded a few minutes later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ATROVENT 500MCG/2ML 20 UNIT DOSE VIALS - row 7
$ws.Range("H7").Value = "1:4"

# DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP. - row 10
$ws.Range("H10").Value = "1:3"

# PULMICORT 0.5MG/ML 20 NEBULIZER VIAL SUSP. - row 13
$ws.Range("H13").Value = "0:14"

# Footer generation timestamp
$ws.Range("A19").Value = "Tuesday, 2 September, 2025 9:54 AM"
